# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Reverses the "Periodo Mora" (period) column for the ERICK ENRIQUE
# VALENZUELA ESCOBAR rows (17-36), and moves the odd due-date value
# (41253) from the old last period row to the new first period row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2311","2310","2309","2308","2307","2306","2305","2304","2303","2302","2301","2212","2211","2210","2209","2208","2207","2206","2205","2204")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 17 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# The special due-date (F column) that used to sit on the last period row
# now belongs to the new first period row, and vice versa.
$ws.Range("F17").Value = 41253
$ws.Range("F36").Value = 44200
